$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.296.15'
$ws.Range("E2").Value = '  +12.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.827.99'
$ws.Range("E3").Value = '  +9.50%  '

$ws.Range("E4").Value = '  -0.52%  '

$ws.Range("E5").Value = '  +4.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.548'
$ws.Range("E6").Value = '  +4.15%  '

$ws.Range("E7").Value = '  -0.46%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.62'
$ws.Range("E8").Value = '  +7.40%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.48'
$ws.Range("E9").Value = '  +7.49%  '

$ws.Range("E10").Value = '  +7.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0671'
$ws.Range("E11").Value = '  +6.02%  '

$ws.Range("E12").Value = '  +2.85%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.088.77'
$ws.Range("E13").Value = '  +9.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.826.21'
$ws.Range("E14").Value = '  +9.69%  '

$ws.Range("E15").Value = '  +5.69%  '

$ws.Range("E16").Value = '  +3.97%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '34.224.57'
$ws.Range("E17").Value = '  +12.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.29'
$ws.Range("E18").Value = '  +8.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.85'
$ws.Range("E19").Value = '  +5.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '259.27'
$ws.Range("E20").Value = '  +7.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0753'
$ws.Range("E21").Value = '  +4.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("E22").Value = '  -0.43%  '

$ws.Range("E23").Value = '  +6.80%  '

$ws.Range("E24").Value = '  +3.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.21'
$ws.Range("E25").Value = '  +3.72%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.48'
$ws.Range("E26").Value = '  +0.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.69'
$ws.Range("E27").Value = '  +5.60%  '

$ws.Range("E28").Value = '  +7.63%  '

$ws.Range("E29").Value = '  +3.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.90'
$ws.Range("E31").Value = '  +12.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0515'
$ws.Range("E32").Value = '  +4.61%  '

$ws.Range("E33").Value = '  +5.84%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.56'
$ws.Range("E34").Value = '  +8.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.551.06'
$ws.Range("E35").Value = '  +3.77%  '

$ws.Range("E36").Value = '  +4.48%  '

$ws.Range("E37").Value = '  +6.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '85.70'
$ws.Range("E38").Value = '  +2.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.631'
$ws.Range("E39").Value = '  +7.11%  '

$ws.Range("E40").Value = '  +5.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.83'
$ws.Range("E41").Value = '  +4.85%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.922'
$ws.Range("E42").Value = '  +10.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.33'
$ws.Range("E43").Value = '  +1.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.17'
$ws.Range("E44").Value = '  +10.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0527'
$ws.Range("E45").Value = '  +5.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.988.56'
$ws.Range("E47").Value = '  +10.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.07'
$ws.Range("E48").Value = '  +161.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.73'
$ws.Range("E49").Value = '  +3.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.18'
$ws.Range("E50").Value = '  +4.17%  '

$ws.Range("E51").Value = '  -0.36%  '
